$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.633.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "'1.847.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").Value = "'312.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").Value = "'0.4270"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("D8").Value = "'0.3624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").Value = "'44.64"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.08%  "

$ws.Range("D10").Value = "'0.07306"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("D11").Value = "'0.8751"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.50%  "

$ws.Range("D12").Value = "'20.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.38%  "

$ws.Range("D13").Value = "'1.876.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "'5.321"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "

$ws.Range("D15").Value = "'6.510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.95%  "

$ws.Range("D16").Value = "'0.06915"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.34%  "

$ws.Range("E17").Value = "  +0.38%  "

$ws.Range("D18").Value = "'79.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.23%  "

$ws.Range("D19").Value = "'0.000009023"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.72%  "

$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "

$ws.Range("E21").Value = "  -0.27%  "

$ws.Range("D22").Value = "'27.653.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").Value = "'4.957"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("E24").Value = "  -3.28%  "

$ws.Range("D25").Value = "'2.127.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.95%  "

$ws.Range("D26").Value = "'1.992"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.64%  "

$ws.Range("D27").Value = "'154.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "

$ws.Range("E28").Value = "  +3.43%  "

$ws.Range("D29").Value = "'122.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.92%  "

$ws.Range("E30").Value = "  -1.07%  "

$ws.Range("D31").Value = "'1.860"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.70%  "

$ws.Range("D32").Value = "'0.08923"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.41%  "

$ws.Range("D33").Value = "'0.7598"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.18%  "

$ws.Range("D34").Value = "'2.966"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.24%  "

$ws.Range("D35").Value = "'4.524"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("E36").Value = "  +1.34%  "

$ws.Range("D37").Value = "'0.05386"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.70%  "

$ws.Range("D38").Value = "'1.090"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").Value = "'0.01933"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.55%  "

$ws.Range("D40").Value = "'2.808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.21%  "

$ws.Range("D41").Value = "'0.5069"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.48%  "

$ws.Range("D42").Value = "'0.1655"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.35%  "

$ws.Range("D43").Value = "'6.752"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").Value = "'8.334"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "

$ws.Range("D45").Value = "'0.06547"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.95%  "

$ws.Range("D46").Value = "'10.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.50%  "

$ws.Range("D47").Value = "'105.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D48").Value = "'0.4669"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("D49").Value = "'1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").Value = "'1.619"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").Value = "'64.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.41%  "
